# Update weekly Fruta/Hortaliza price data (Mora) rows 2-14.
# Only columns D (Fecha), M (Volumen), N (Precio minimo), O (Precio maximo),
# P (Precio promedio ponderado), R (Origen) and S (Precio $/Kg) change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param(
        [int]$RowNum,
        [double]$DVal,
        [double]$MVal,
        [double]$NVal,
        [double]$OVal,
        [double]$PVal,
        [string]$RVal,
        [double]$SVal
    )

    $ws.Range("D$RowNum").Value = $DVal
    $ws.Range("M$RowNum").Value = $MVal
    $ws.Range("N$RowNum").Value = $NVal
    $ws.Range("O$RowNum").Value = $OVal
    $ws.Range("P$RowNum").Value = $PVal
    $ws.Range("R$RowNum").Value = $RVal
    $ws.Range("S$RowNum").Value = $SVal
}

Set-Row 2  44617 90  6500 6500 6500 "Provincia de Curicó"  3250
Set-Row 3  44208 85  3000 3000 3000 "Provincia de Linares" 1500
Set-Row 4  44188 150 3000 3400 3240 "Provincia de Linares" 1620
Set-Row 5  44236 300 3600 4000 3800 "Provincia de Curicó"  1900
Set-Row 6  44586 250 5000 5000 5000 "Provincia de Curicó"  2500
Set-Row 7  44168 170 8000 8000 8000 "Provincia de Linares" 4000
Set-Row 9  44582 380 5000 5000 5000 "Provincia de Curicó"  2500
Set-Row 11 44194 120 3000 3000 3000 "Provincia de Linares" 1500
Set-Row 12 44174 200 3200 3200 3200 "Provincia de Curicó"  1600
Set-Row 13 44237 100 3600 4000 3800 "Provincia de Curicó"  1900
Set-Row 14 44231 150 3400 3400 3400 "Provincia de Curicó"  1700
